$wb = $excel.ActiveWorkbook

# --- Sheet1: CombinationsWithRepeat -----------------------------------
$ws1 = $wb.Worksheets.Item("CombinationsWithRepeat")
$ws1.Activate()

# New permutations-with-repetition example (rows 18-19)
$ws1.Range("D18").Value = "n="
$ws1.Range("E18").Value = 26
$ws1.Range("D19").Value = "r="
$ws1.Range("E19").Value = 24
$ws1.Range("H18").Formula = "=(FACT(E19+E18-1))/(FACT(+E19)*(FACT(E18-1)))"

# Column H is widened to fit the formula text (32.28515625 OOXML units,
# expressed as a COM "characters" width so the host's internal MDW-6
# pixel-grid round-trips back to that exact stored width)
$ws1.Columns.Item(8).ColumnWidth = 31.451822916666668

# H18 displays as a whole number (format code "0")
$ws1.Range("H18").NumberFormat = "0"

# Update the selection left behind on this sheet
$ws1.Range("H18").Select() | Out-Null

# --- Sheet2: CombinationsWITHOUTRepeat ---------------------------------
$ws2 = $wb.Worksheets.Item("CombinationsWITHOUTRepeat")
$ws2.Activate()
$ws2.Range("C14").Select() | Out-Null

# --- Sheet3: permutationsWITHRepeat (new, blank, becomes active tab) --
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "permutationsWITHRepeat"
$ws3.Select()
